$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header texts in row 10 to append " Rs." to the vendor basic charge
# and vendor tax column headers.
$ws.Range("G10").Value = "Vendor " + [char]10 + "Basic Charge Rs."
$ws.Range("H10").Value = "Vendor " + [char]10 + "Tax Rs."

# Move the active selection to G10 (matches the workbook's saved cursor position).
$ws.Range("G10").Select()
